# Improve robustness of data loader: the "Branch" column (column C) in the
# "Template" sheet is no longer used by the loader, so remove it entirely.
# This shifts every subsequent header cell one column to the left
# (D->C, E->D, ... T->S) and drops the sheet's used range from A1:T1 to
# A1:S1. The shared string "Branch" itself is left in sharedStrings.xml
# because it is still referenced by the "Sample Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Template")

$ws.Columns("C:C").Delete()
